$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the longer task names
$ws.Columns.Item(1).ColumnWidth = 36.5703125

# Update row 4 (Modelling CNN with overfitting handling)
$ws.Range("A4").Value = "Modelling CNN with overfiting handiling"
$ws.Range("B4").Value = Get-Date -Year 2021 -Month 11 -Day 8 -Hour 0 -Minute 0 -Second 0
$ws.Range("C4").Value = Get-Date -Year 2021 -Month 11 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Range("D4").Value = "Jody"
$ws.Range("E4").Value = "Done"

# Update row 5 (Model Deployment)
$ws.Range("A5").Value = "Model Deploymet"
$ws.Range("B5").Value = Get-Date -Year 2021 -Month 11 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Range("C5").Value = Get-Date -Year 2021 -Month 11 -Day 12 -Hour 0 -Minute 0 -Second 0
$ws.Range("D5").Value = "Aries"
$ws.Range("E5").Value = "On Going"

# Update statuses on rows 2 & 3 to "Done"
$ws.Range("E2").Value = "Done"
$ws.Range("E3").Value = "Done"

# Add two new empty rows (6 and 7) matching row 5's formatting (no values)
$ws.Range("A5:E5").Copy()
$ws.Range("A6:E7").PasteSpecial(-4122)

$ws.Range("C8").Select()
